# Update the "views/attendance" counters in column F across the three
# event-listing sheets (展览, 演出, 全部类型). 本地生活 has no F data and is untouched.

$wb = $excel.ActiveWorkbook

# --- 展览 (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 10341
$ws1.Range("F9").Value  = 775
$ws1.Range("F12").Value = 1088
$ws1.Range("F13").Value = 3225
$ws1.Range("F33").Value = 10
$ws1.Range("F38").Value = 465
$ws1.Range("F39").Value = 449
$ws1.Range("F43").Value = 52
$ws1.Range("F44").Value = 456

# --- 演出 (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 47

# --- 全部类型 (All types, combined listing) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 10341
$ws4.Range("F9").Value  = 775
$ws4.Range("F10").Value = 1088
$ws4.Range("F11").Value = 3225
$ws4.Range("F28").Value = 10
$ws4.Range("F29").Value = 47
$ws4.Range("F38").Value = 449
$ws4.Range("F45").Value = 52
$ws4.Range("F46").Value = 456
